$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 116
$ws.Range("B116").Value = 6798562
$ws.Range("F116").Value = "AC Horsens"
$ws.Range("G116").Value = "Kolding IF"
$ws.Range("H116").Value = 1
$ws.Range("I116").Value = 2
$ws.Range("J116").Value = "A"
$ws.Range("K116").Value = 2.4
$ws.Range("L116").Value = 3.5
$ws.Range("M116").Value = 2.4
$ws.Range("N116").Value = 3.4
$ws.Range("O116").Value = 3.3
$ws.Range("P116").Value = 2.15
$ws.Range("Q116").Value = 0.25
$ws.Range("R116").Value = 1.975
$ws.Range("S116").Value = 1.875
$ws.Range("U116").Value = 1.975
$ws.Range("V116").Value = 1.875
$ws.Range("X116").Value = -1
$ws.Range("Y116").Value = 1.15
$ws.Range("Z116").Value = -1
$ws.Range("AA116").Value = 0.875
$ws.Range("AB116").Value = 0.9750000000000001

# Row 117
$ws.Range("B117").Value = 6799307
$ws.Range("F117").Value = "Vendsyssel FF"
$ws.Range("G117").Value = "Hillerd"
$ws.Range("H117").Value = 4
$ws.Range("I117").Value = 4
$ws.Range("J117").Value = "D"
$ws.Range("K117").Value = 1.75
$ws.Range("L117").Value = 3.6
$ws.Range("M117").Value = 4.333
$ws.Range("N117").Value = 2.4
$ws.Range("O117").Value = 3.2
$ws.Range("P117").Value = 3
$ws.Range("Q117").Value = -0.25
$ws.Range("R117").Value = 2.1
$ws.Range("S117").Value = 1.775
$ws.Range("U117").Value = 1.85
$ws.Range("V117").Value = 2
$ws.Range("X117").Value = 2.2
$ws.Range("Y117").Value = -1
$ws.Range("Z117").Value = -0.5
$ws.Range("AA117").Value = 0.3875
$ws.Range("AB117").Value = 0.8500000000000001

# Row 135
$ws.Range("B135").Value = 6799316
$ws.Range("F135").Value = "Vendsyssel FF"
$ws.Range("G135").Value = "Hobro IK"
$ws.Range("H135").Value = 1
$ws.Range("I135").Value = 1
$ws.Range("K135").Value = 1.909
$ws.Range("L135").Value = 3.4
$ws.Range("M135").Value = 3.75
$ws.Range("N135").Value = 1.909
$ws.Range("O135").Value = 3.5
$ws.Range("P135").Value = 4.2
$ws.Range("Q135").Value = -0.5
$ws.Range("R135").Value = 1.925
$ws.Range("S135").Value = 1.925
$ws.Range("U135").Value = 2
$ws.Range("V135").Value = 1.85
$ws.Range("X135").Value = 2.5
$ws.Range("AA135").Value = 0.925
$ws.Range("AB135").Value = -0.5
$ws.Range("AC135").Value = 0.425

# Row 136
$ws.Range("B136").Value = 6798566
$ws.Range("F136").Value = "Kolding IF"
$ws.Range("G136").Value = "Naestved"
$ws.Range("H136").Value = 3
$ws.Range("I136").Value = 3
$ws.Range("J136").Value = "D"
$ws.Range("K136").Value = 1.666
$ws.Range("M136").Value = 4.333
$ws.Range("N136").Value = 1.7
$ws.Range("P136").Value = 5.25
$ws.Range("Q136").Value = -0.75
$ws.Range("R136").Value = 1.975
$ws.Range("S136").Value = 1.875
$ws.Range("T136").Value = 2.25
$ws.Range("U136").Value = 1.8
$ws.Range("V136").Value = 2.05
$ws.Range("X136").Value = 2.6
$ws.Range("Y136").Value = -1
$ws.Range("AA136").Value = 0.875
$ws.Range("AB136").Value = 0.8

# Row 137
$ws.Range("B137").Value = 6798882
$ws.Range("F137").Value = "Hillerd"
$ws.Range("G137").Value = "AaB"
$ws.Range("H137").Value = 1
$ws.Range("J137").Value = "A"
$ws.Range("K137").Value = 4
$ws.Range("L137").Value = 3.8
$ws.Range("M137").Value = 1.75
$ws.Range("N137").Value = 3.6
$ws.Range("O137").Value = 3.6
$ws.Range("P137").Value = 2
$ws.Range("Q137").Value = 0.5
$ws.Range("R137").Value = 1.825
$ws.Range("S137").Value = 2.025
$ws.Range("T137").Value = 2.5
$ws.Range("U137").Value = 1.825
$ws.Range("V137").Value = 2.025
$ws.Range("X137").Value = -1
$ws.Range("Y137").Value = 1
$ws.Range("AA137").Value = 1.025
$ws.Range("AB137").Value = 0.825

# Row 138
$ws.Range("B138").Value = 6799314
$ws.Range("F138").Value = "FC Fredericia"
$ws.Range("G138").Value = "HB Kge"
$ws.Range("I138").Value = 2
$ws.Range("J138").Value = "D"
$ws.Range("K138").Value = 1.444
$ws.Range("L138").Value = 4.75
$ws.Range("M138").Value = 5.5
$ws.Range("N138").Value = 1.5
$ws.Range("O138").Value = 4.75
$ws.Range("P138").Value = 5.5
$ws.Range("Q138").Value = -1.25
$ws.Range("R138").Value = 1.975
$ws.Range("S138").Value = 1.875
$ws.Range("W138").Value = -1
$ws.Range("X138").Value = 3.75
$ws.Range("Z138").Value = -1
$ws.Range("AA138").Value = 0.875
$ws.Range("AB138").Value = 1.025
$ws.Range("AC138").Value = -1

# Row 139
$ws.Range("B139").Value = 6799315
$ws.Range("F139").Value = "Sonderjyske"
$ws.Range("G139").Value = "FC Helsingor"
$ws.Range("H139").Value = 2
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = "H"
$ws.Range("K139").Value = 1.3
$ws.Range("L139").Value = 5.5
$ws.Range("M139").Value = 7
$ws.Range("N139").Value = 1.3
$ws.Range("O139").Value = 5.75
$ws.Range("P139").Value = 8.5
$ws.Range("Q139").Value = -1.5
$ws.Range("R139").Value = 1.85
$ws.Range("S139").Value = 2
$ws.Range("T139").Value = 3.25
$ws.Range("U139").Value = 2.025
$ws.Range("V139").Value = 1.825
$ws.Range("W139").Value = 0.3
$ws.Range("X139").Value = -1
$ws.Range("Z139").Value = 0.8500000000000001
$ws.Range("AA139").Value = -1
$ws.Range("AB139").Value = -1
$ws.Range("AC139").Value = 0.825

# Row 146
$ws.Range("R146").Value = 1.875
$ws.Range("S146").Value = 1.975
$ws.Range("T146").Value = 3
$ws.Range("U146").Value = 2.05
$ws.Range("V146").Value = 1.8

# Row 147
$ws.Range("P147").Value = 2.15
$ws.Range("T147").Value = 2.75
$ws.Range("U147").Value = 2.025
$ws.Range("V147").Value = 1.825

# Row 148
$ws.Range("N148").Value = 2.2
$ws.Range("P148").Value = 3.25
$ws.Range("R148").Value = 1.9
$ws.Range("S148").Value = 1.95
$ws.Range("U148").Value = 1.975
$ws.Range("V148").Value = 1.875

# Row 149
$ws.Range("O149").Value = 4.75
$ws.Range("P149").Value = 7
$ws.Range("R149").Value = 2.05
$ws.Range("S149").Value = 1.8
$ws.Range("T149").Value = 2.5
$ws.Range("U149").Value = 1.825
$ws.Range("V149").Value = 2.025

# Row 150
$ws.Range("P150").Value = 4.2

# Row 151 (new row - copy formatting from row 150, then set values)
$ws.Range("A150:AC150").Copy()
$ws.Range("A151:AC151").PasteSpecial(-4122)
$ws.Range("H151").ClearContents()
$ws.Range("I151").ClearContents()
$ws.Range("J151").ClearContents()
$ws.Range("AB151").ClearContents()
$ws.Range("AC151").ClearContents()

$ws.Range("A151").Value = 149
$ws.Range("B151").Value = 7993177
$ws.Range("C151").Value = "Denmark Division 1"
$ws.Range("D151").Value = "Denmark Division 1"
$ws.Range("E151").Value = 45389.41666666666
$ws.Range("F151").Value = "Vendsyssel FF"
$ws.Range("G151").Value = "AaB"
$ws.Range("K151").Value = 2.8
$ws.Range("L151").Value = 3.4
$ws.Range("M151").Value = 2.375
$ws.Range("N151").Value = 3.5
$ws.Range("O151").Value = 3.5
$ws.Range("P151").Value = 2.1
$ws.Range("Q151").Value = 0.25
$ws.Range("R151").Value = 2.025
$ws.Range("S151").Value = 1.825
$ws.Range("T151").Value = 2.5
$ws.Range("U151").Value = 1.825
$ws.Range("V151").Value = 2.025
$ws.Range("W151").Value = 0
$ws.Range("X151").Value = 0
$ws.Range("Y151").Value = 0
$ws.Range("Z151").Value = 0
$ws.Range("AA151").Value = 0
